# Update countries & provincias Spain
# - Refresh the "last updated" timestamp
# - Refresh totals for Estados Unidos (row 4) and Moldavia (row 59)
# - Reorder a few small countries (swap Montserrat/Islas Malvinas, move
#   Cabo Verde up in the list) updating their stats, and refresh
#   "Provincias Espana" rows 193-198

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 14 de Abril de 2020 a las 22:22"

# --- Estados Unidos (row 4) -------------------------------------------
$ws.Range("B4").Value = 608988
$ws.Range("C4").Value = 22047
$ws.Range("D4").Value = 38479
$ws.Range("E4").Value = 544838
$ws.Range("F4").Value = 12829
$ws.Range("G4").Value = 2031
$ws.Range("H4").Value = 25671

# --- Moldavia (row 59) --------------------------------------------------
$ws.Range("E59").Value = 1760
$ws.Range("G59").Value = 5
$ws.Range("H59").Value = 40

# --- Reorder / update small countries (rows 193-198) --------------------
# Row 193: was "Islas Malvinas" -> now "Montserrat" (unchanged stats)
$ws.Range("A193").Value = "Montserrat"
$ws.Range("B193").Value = 11
$ws.Range("C193").Value = 0
$ws.Range("D193").Value = 1
$ws.Range("E193").Value = 10
$ws.Range("F193").Value = 1
$ws.Range("G193").Value = 0
$ws.Range("H193").Value = 0

# Row 194: was "Montserrat" -> now "Islas Malvinas" (unchanged stats)
$ws.Range("A194").Value = "Islas Malvinas"
$ws.Range("B194").Value = 11
$ws.Range("C194").Value = 6
$ws.Range("D194").Value = 1
$ws.Range("E194").Value = 10
$ws.Range("F194").Value = 0
$ws.Range("G194").Value = 0
$ws.Range("H194").Value = 0

# Row 195: was "Republica de Africa Central" -> now "Cabo Verde" (updated stats)
$ws.Range("A195").Value = "Cabo Verde"
$ws.Range("B195").Value = 11
$ws.Range("C195").Value = 1
$ws.Range("D195").Value = 1
$ws.Range("E195").Value = 9
$ws.Range("F195").Value = 0
$ws.Range("G195").Value = 0
$ws.Range("H195").Value = 1

# Row 196: was "Groenlandia" -> now "Republica de Africa Central" (unchanged stats, shifted)
$ws.Range("A196").Value = "Republica de Africa Central"
$ws.Range("B196").Value = 11
$ws.Range("C196").Value = 0
$ws.Range("D196").Value = 4
$ws.Range("E196").Value = 7
$ws.Range("F196").Value = 0
$ws.Range("G196").Value = 0
$ws.Range("H196").Value = 0

# Row 197: was "Islas Turcas y Caicos" -> now "Groenlandia" (unchanged stats, shifted)
$ws.Range("A197").Value = "Groenlandia"
$ws.Range("B197").Value = 11
$ws.Range("C197").Value = 0
$ws.Range("D197").Value = 11
$ws.Range("E197").Value = 0
$ws.Range("F197").Value = 0
$ws.Range("G197").Value = 0
$ws.Range("H197").Value = 0

# Row 198: was "Cabo Verde" -> now "Islas Turcas y Caicos" (unchanged stats, shifted)
$ws.Range("A198").Value = "Islas Turcas y Caicos"
$ws.Range("B198").Value = 10
$ws.Range("C198").Value = 0
$ws.Range("D198").Value = 0
$ws.Range("E198").Value = 9
$ws.Range("F198").Value = 0
$ws.Range("G198").Value = 0
$ws.Range("H198").Value = 1
